$wb = $excel.ActiveWorkbook

$wsWeather = $wb.Worksheets.Item("WeatherData")
$wsWeather.Columns.Item(5).Delete()

$wsBeachLog = $wb.Worksheets.Item("BeachLog")
$wsBeachLog.Columns.Item(2).Delete()
